# Fruta / hortaliza, semanal
# Shifts rows 122-154 down by one (for the weekly price columns) and
# appends a new row 155 that carries what used to be row 154's data.
# Row 122 itself receives a brand-new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "O", "P")

# 1) Snapshot the ORIGINAL values for the columns that move, for rows 121-154,
#    before any writes happen (so later writes never clobber a value we still
#    need to read).
$orig = @{}
for ($r = 121; $r -le 154; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range($c + $r).Value()
    }
    $orig[$r] = $row
}

# 2) Append row 155 = original row 154 (every column, not just the movers).
$ws.Range("A155").Value = $ws.Range("A154").Value()
$ws.Range("B155").Value = $ws.Range("B154").Value()
$ws.Range("C155").Value = $ws.Range("C154").Value()
$ws.Range("D155").Value = $orig[154]["D"]
$ws.Range("D155").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E155").Value = $ws.Range("E154").Value()
$ws.Range("F155").Value = $ws.Range("F154").Value()
$ws.Range("G155").Value = $ws.Range("G154").Value()
$ws.Range("H155").Value = $ws.Range("H154").Value()
$ws.Range("I155").Value = $ws.Range("I154").Value()
$ws.Range("J155").Value = $orig[154]["J"]
$ws.Range("K155").Value = $orig[154]["K"]
$ws.Range("L155").Value = $orig[154]["L"]
$ws.Range("M155").Value = $orig[154]["M"]
$ws.Range("N155").Value = $ws.Range("N154").Value()
$ws.Range("O155").Value = $orig[154]["O"]
$ws.Range("P155").Value = $orig[154]["P"]
$ws.Range("Q155").Value = $ws.Range("Q154").Value()
$ws.Range("R155").Value = $ws.Range("R154").Value()

# 3) Shift: row N (123..154) takes the ORIGINAL row (N-1) values for the
#    moving columns.
for ($r = 154; $r -ge 123; $r--) {
    $prev = $orig[$r - 1]
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value = $prev[$c]
    }
}

# 4) Row 122 gets the brand-new weekly observation (O122 keeps its existing
#    value, "Región Metropolitana" - unchanged by the diff).
$ws.Range("D122").Value = 44551
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 2500
$ws.Range("L122").Value = 2500
$ws.Range("M122").Value = 2500
$ws.Range("P122").Value = 833
